$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the title/header text cells (E2, E3:E6, E7, E8)
$ws.Range("E2").ClearContents()
$ws.Range("E3:E6").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("E8").ClearContents()

# Clear the category-mapping values in column H (rows 11-54)
$ws.Range("H11:H54").ClearContents()

# Clear the check formula in H56
$ws.Range("H56").ClearContents()

# Update the selected range shown in the sheet view to the whole column H
$ws.Columns("H").Select()
